$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01178416625306948
$ws.Range("C2").Value = 0.1962658956258756

$ws.Range("B3").Value = 0.09240034103935675
$ws.Range("C3").Value = 0.2022908011725063

$ws.Range("B4").Value = 0.675461405856044
$ws.Range("C4").Value = 0.1324541629393382

$ws.Range("B5").Value = 0.9729398810865978
$ws.Range("C5").Value = 0.4335099445219128

$ws.Range("B6").Value = 0.8494076316771882
$ws.Range("C6").Value = 0.5920963804129344

$ws.Range("B7").Value = 0.7088019934232593
$ws.Range("C7").Value = 0.07778007816758685

$ws.Range("B8").Value = 0.004472704529762268
$ws.Range("C8").Value = 0.2338663482666016

$ws.Range("B9").Value = 0.07509490836423421
$ws.Range("C9").Value = 0.1591592248205764

$ws.Range("B10").Value = 0.7669162226934052
$ws.Range("C10").Value = 0.5177064802775271
